# Commit: "skip data test in one testcase"
#
# - login_003 sheet: Runmode (D2) flips from "Y" to "N" so the test runner
#   skips the data-driven test in this testcase.
# - Selection/active-sheet bookkeeping follows the author's last clicks:
#   the "login_003" sheet view loses its tabSelected flag, and the
#   "Test Steps" sheet becomes the active/selected tab instead.

$wb = $excel.ActiveWorkbook

# login_003 sheet: flip the Runmode flag from Y to N, then leave the
# selection on C2 (where the author's cursor ended up).
$wsLogin = $wb.Worksheets.Item("login_003")
$wsLogin.Range("D2").Value = "N"
$wsLogin.Range("C2").Select()

# Make "Test Steps" the active sheet/tab, with C6 selected.
$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsSteps.Activate()
$wsSteps.Range("C6").Select()
